$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: "time_taken" - style matches the other header cells (copy format from E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Timestamps for each data row F2:F16 (plain, unstyled cells like the rest of the data rows)
$timestamps = @(
    "2021-10-05 10:50:08.641891",
    "2021-10-05 10:50:08.641901",
    "2021-10-05 10:50:08.641904",
    "2021-10-05 10:50:08.641907",
    "2021-10-05 10:50:08.641910",
    "2021-10-05 10:50:08.641913",
    "2021-10-05 10:50:08.641916",
    "2021-10-05 10:50:08.641919",
    "2021-10-05 10:50:08.641922",
    "2021-10-05 10:50:08.641924",
    "2021-10-05 10:50:08.641927",
    "2021-10-05 10:50:08.641929",
    "2021-10-05 10:50:08.641932",
    "2021-10-05 10:50:08.641935",
    "2021-10-05 10:50:08.641937"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Output "applied time_taken column"
